$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new columns before column W (they become the new W and X columns),
# pushing the old VIN..AdditionalDriver block from W..AE to Y..AG.
$ws.Range("W1:X1").EntireColumn.Insert()

# Fill the new SSN / DL columns, matching the order the author appears to have
# used (W column top-to-bottom, then X column top-to-bottom for the data rows).
$ws.Range("W1").Value = "SSN"
$ws.Range("X1").Value = "DL"

$ws.Range("W2").Value = "666195140"
$ws.Range("X2").Value = "487956891"

$ws.Range("W3").Value = "666195141"
$ws.Range("W4").Value = "666195142"
$ws.Range("W5").Value = "666195143"

$ws.Range("X3").Value = "487956892"
$ws.Range("X4").Value = "487956893"
$ws.Range("X5").Value = "487956894"

$ws.Range("D6").Value = "htcglobal2019"

$ws.Range("W6").Value = "666195144"
$ws.Range("X6").Value = "487956895"

# Match styling used by the other data cells in these columns (Consolas text format).
$ws.Range("W2:X6").Style = $ws.Range("G2").Style

# Row 3's explicit row height is cleared back to the sheet default.
$ws.Rows.Item(3).RowHeight = $ws.Rows.Item(4).RowHeight

# Sheet view bookkeeping to match the saved workbook state.
$ws.Application.ActiveWindow.ScrollColumn = 11
$ws.Range("X13").Select()
